# Convert the year header row (B1:Q1) from text labels ("A2010" ... "A2025")
# into plain numeric years (2010 ... 2025), and align the style/format of the
# last two header cells (P1, Q1) with the rest of the header row (B1:O1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$year = 2010
foreach ($col in 2..17) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $year
    $year++
}

# P1 and Q1 previously used a slightly different style/format than B1:O1;
# copy the formatting from B1 onto them so the whole header row matches.
$ws.Range("B1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to H16, matching the author's final cursor spot.
$ws.Range("H16").Select()
